# First implementation for finding the hazard rate
# - Keep the existing (formula-driven) "Curva" sheet but rename it to
#   "Curva anterior" (the previous curve).
# - Insert a brand-new "Curva" sheet in its place, populated with a
#   static snapshot of the current zero-rate curve (dates + ZR values),
#   styled with a bold bordered header and bordered/Arial data cells.
# - Tidy up a couple of incidental selection changes left behind on the
#   "Aux" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the current "Curva" (formula) sheet to "Curva anterior" and
#    add a fresh "Curva" sheet right in front of it.
# ---------------------------------------------------------------------
$wsOldCurva = $wb.Worksheets.Item("Curva")
$wsOldCurva.Name = "Curva anterior"

$wsCurva = $wb.Worksheets.Add([System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 1, [System.Reflection.Missing]::Value)
$wsCurva.Name = "Curva"
$wsCurva.Move([System.Reflection.Missing]::Value, $wsOldCurva)

# ---------------------------------------------------------------------
# 2. Fix up the "Curva anterior" sheet view (it keeps the old data /
#    formulas untouched, just the selection + active-tab state moves).
# ---------------------------------------------------------------------
$wsOldCurva.Range("A1:B1").Select() | Out-Null

$wsAux = $wb.Worksheets.Item("Aux")
$wsAux.Range("B3").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Populate the new "Curva" sheet with the static curve snapshot.
# ---------------------------------------------------------------------
$wsCurva.Columns.Item(1).ColumnWidth = 15.5

$wsCurva.Range("A1").Value2 = "Fecha"
$wsCurva.Range("B1").Value2 = "ZR"

$data = @(
    @(2, 45954, 3.9891399999999999),
    @(3, 45965, 3.9891399999999999),
    @(4, 45973, 3.9596),
    @(5, 45979, 3.9491900000000002),
    @(6, 45989, 3.9392299999999998),
    @(7, 46020, 3.8572199999999999),
    @(8, 46050, 3.7948400000000002),
    @(9, 46080, 3.7295699999999998),
    @(10, 46111, 3.6789200000000002),
    @(11, 46140, 3.6342400000000001),
    @(12, 46170, 3.5933000000000002),
    @(13, 46202, 3.5539200000000002),
    @(14, 46231, 3.5169600000000001),
    @(15, 46262, 3.47776),
    @(16, 46293, 3.4430800000000001),
    @(17, 46323, 3.41066),
    @(18, 46505, 3.2620200000000001),
    @(19, 46688, 3.1904300000000001),
    @(20, 47056, 3.1570299999999998),
    @(21, 47420, 3.1829900000000002),
    @(22, 47784, 3.2283400000000002),
    @(23, 48149, 3.2840099999999999),
    @(24, 48515, 3.3460200000000002),
    @(25, 48880, 3.40883),
    @(26, 49247, 3.47241),
    @(27, 49611, 3.5341),
    @(28, 50341, 3.65178),
    @(29, 51438, 3.7943600000000002),
    @(30, 53265, 3.9108800000000001),
    @(31, 55089, 3.9050600000000002),
    @(32, 56915, 3.8376899999999998),
    @(33, 60568, 3.6242100000000002),
    @(34, 64220, 3.3721199999999998)
)

foreach ($item in $data) {
    $r = $item[0]
    $wsCurva.Cells.Item($r, 1).Value2 = $item[1]
    $wsCurva.Cells.Item($r, 2).Value2 = $item[2]
}

# Date formatting for column A (built-in numFmtId 14).
$wsCurva.Range("A2:A34").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------
# 4. Styling: bold bordered header row, bordered Arial data rows, with
#    the first data row (row 2) kept bold to match the source workbook.
# ---------------------------------------------------------------------
$full = $wsCurva.Range("A1:B34")
$full.Font.Name = "Arial"
$full.VerticalAlignment = -4108   # xlCenter
$full.WrapText = $true
$full.Borders.Weight = -4138      # xlMedium
$full.Borders.Color = 0
$full.RowHeight = 15

$header = $wsCurva.Range("A1:B1")
$header.Font.Bold = $true

$row2Bold = $wsCurva.Range("A2")
$row2Bold.Font.Bold = $true
$row2Bold.Font.Name = "Arial"

$wsCurva.Range("A2:A34").NumberFormat = "mm-dd-yy"

$wsCurva.Range("A2").Select() | Out-Null
$wsCurva.Activate()
